# ---------------------------------------------------------------------------
# Applies the edit described by the diff:
#  - In the last "ListParagraph" bullet (about Impostor Pass Rate / False
#    Reject Rates), remove the italic formatting and rewrite the trailing
#    sentence.
#  - Insert a new Heading1 paragraph "False Accept Rate".
#  - Insert a new body paragraph explaining the False Accept Rate, written
#    as two separate runs, and move the "_GoBack" bookmark onto the end of
#    this new paragraph.
#  - Drop one of the two trailing empty paragraphs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- locate the target bullet paragraph (the "Finally, will be taking..." one) ---
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("Finally, will be taking")) {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq 0) { throw "target paragraph not found" }

# --- 1) rewrite the two runs of that paragraph ---
$oldRun1 = "Finally, will be taking a look at the various genuine scores and using a threshold T and using the same threshold T, and looking at the various impostor scores, will calculate the Impostor Pass Rate and False Reject Rates. "
$newRun1 = "Finally, will be taking a look at the various genuine scores and using a threshold T and using the same threshold T, and looking at the various impostor scores, will calculate the Impostor Pass Rate and False Reject Rates"
$rng = $d.Paragraphs($targetIdx).Range
$rng.Find.Execute($oldRun1, $false, $false, $false, $false, $false, $true, 1, $false, $newRun1, 2) | Out-Null

$oldRun2 = " – This step I have tried, but due to time constraints I was not able to accomplish this step.  "
$newRun2 = ".  – This step, please note I will only be able have it such that I can find the Impostor Pass Rate which can also be referred to as the False Accept Rate.  I will briefly talk about the concept of this measure "
$rng2 = $d.Paragraphs($targetIdx).Range
$rng2.Find.Execute($oldRun2, $false, $false, $false, $false, $false, $true, 1, $false, $newRun2, 2) | Out-Null

# --- 2) drop the italics on that whole paragraph (covers both runs and the
#        paragraph mark's own run properties) ---
$d.Paragraphs($targetIdx).Range.Font.Italic = $false

# --- 3) insert the new Heading1 paragraph right after it ---
$d.Paragraphs($targetIdx).Range.InsertParagraphAfter()
$headIdx = $targetIdx + 1
$d.Paragraphs($headIdx).Range.Text = "False Accept Rate"
$d.Paragraphs($headIdx).Style = "Normal"
$d.Paragraphs($headIdx).Style = "Heading1"

# --- 4) build the new explanatory paragraph by duplicating a plain
#        (Normal-style, two plain runs) paragraph already in the document,
#        then overwriting its text - this keeps the "two separate <w:r>"
#        shape instead of collapsing to a single run. ---
$templateIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("For this programming assignment")) {
        $templateIdx = $i
        break
    }
}
if ($templateIdx -eq 0) { throw "template paragraph not found" }

$tStart = $d.Paragraphs($templateIdx).Range.Start
$tEnd = $d.Paragraphs($templateIdx + 1).Range.Start
$d.Range($tStart, $tEnd).Copy()

$headEnd = $d.Paragraphs($headIdx).Range.End
$d.Range($headEnd, $headEnd).Paste()
$bodyIdx = $headIdx + 1

$oldBody1 = "For this programming assignment, I have used the C#.net programming language, the Visual Studio 2013 Ultimate IDE for Windows 8.1 along with the PasswordData.csv file which I manipulated into 51 separate csv files, and when running the program, I have the code written such that depending upon the user number and the value of N, which represents the sample size: the training and testing vectors are extracted from the original data, the template vectors, genuine scores, impostor scores, and the various rates are calculated appropriately. "
$newBody1 = "This is sometimes referred to as the impostor pass rate, to find this measure, I will be taking a threshold that the end user will put in and in all of the impostor scores that I have calculated, I will compare the value of the user entered threshold and the score.  If the impostor score is less than or equal to the threshold value, "
$oldBody2 = "The approach that I have taken for this assignment is detailed below:"
$newBody2 = "that means the impostor has passed through the system. "

$bodyStart = $d.Paragraphs($bodyIdx).Range.Start
$sub1 = $d.Range($bodyStart, $bodyStart + $oldBody1.Length)
$sub1.Text = $newBody1

# put a transient bookmark at the boundary between the two new runs so the
# engine does not silently re-merge them into a single run once both share
# the same (empty) run formatting
$boundary = $bodyStart + $newBody1.Length
$d.Bookmarks.Add("ZZTMP", $d.Range($boundary, $boundary)) | Out-Null

$sub2 = $d.Range($boundary, $boundary + $oldBody2.Length)
$sub2.Text = $newBody2

$d.Bookmarks("ZZTMP").Delete()

# --- 5) move the "_GoBack" bookmark from the old paragraph onto the end of
#        the new explanatory paragraph ---
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$bodyEnd = $d.Paragraphs($bodyIdx).Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bodyEnd, $bodyEnd)) | Out-Null

# --- 6) remove one of the two trailing empty paragraphs ---
$lastCount = $d.Paragraphs.Count
$d.Paragraphs($lastCount).Range.Delete()

Write-Host "done"
